$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 905.95654
$ws.Range("I92").Value = 849.381
$ws.Range("K92").Value = 849.381
$ws.Range("M92").Value = 398.619
$ws.Range("H100").Value = 1971.4
$ws.Range("I100").Value = 1401.875
$ws.Range("K100").Value = 1401.875
$ws.Range("M100").Value = -860.875
$ws.Range("H116").Value = 10657.177
$ws.Range("I116").Value = 9430
$ws.Range("J116").Value = 13602.4
$ws.Range("K116").Value = 9430
$ws.Range("L116").Value = 13602.4
$ws.Range("M116").Value = -5988
$ws.Range("N116").Value = -20486.4
$ws.Range("H138").Value = 4354.9
$ws.Range("I138").Value = 4314.846
$ws.Range("J138").Value = 4385.5293
$ws.Range("K138").Value = 12944.538
$ws.Range("L138").Value = 13156.5879
$ws.Range("M138").Value = -7804.537999999999
$ws.Range("N138").Value = -23436.5879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1740.5834
$ws.Range("I2").Value = 2349.8333
$ws.Range("K2").Value = 2349.8333
$ws.Range("M2").Value = -2236.8333
$ws.Range("H97").Value = 2330.5
$ws.Range("I97").Value = 602.375
$ws.Range("J97").Value = 5786.75
$ws.Range("K97").Value = 602.375
$ws.Range("L97").Value = 5786.75
$ws.Range("M97").Value = -106.375
$ws.Range("N97").Value = -6778.75
$ws.Range("H102").Value = 1387.8667
$ws.Range("I102").Value = 1455.5385
$ws.Range("J102").Value = 948
$ws.Range("K102").Value = 1455.5385
$ws.Range("L102").Value = 948
$ws.Range("M102").Value = 166.4614999999999
$ws.Range("N102").Value = -4192
$ws.Range("H116").Value = 1740.5834
$ws.Range("I116").Value = 2349.8333
$ws.Range("K116").Value = 2349.8333
$ws.Range("M116").Value = -55.83329999999978

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1740.5834
$ws.Range("I3").Value = 2349.8333
$ws.Range("K3").Value = 2349.8333
$ws.Range("M3").Value = -2235.8333
$ws.Range("H86").Value = 27032.691
$ws.Range("I86").Value = 24662.6
$ws.Range("J86").Value = 34933
$ws.Range("K86").Value = 24662.6
$ws.Range("L86").Value = 34933
$ws.Range("M86").Value = -23539.6
$ws.Range("N86").Value = -37179
$ws.Range("H89").Value = 27032.691
$ws.Range("I89").Value = 24662.6
$ws.Range("J89").Value = 34933
$ws.Range("K89").Value = 123313
$ws.Range("L89").Value = 174665
$ws.Range("M89").Value = -117697
$ws.Range("N89").Value = -185897
$ws.Range("H107").Value = 2784.2856
$ws.Range("I107").Value = 2797.2
$ws.Range("K107").Value = 2797.2
$ws.Range("M107").Value = -877.1999999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 31649.143
$ws.Range("J41").Value = 35990.668
$ws.Range("L41").Value = 35990.668
$ws.Range("N41").Value = -36846.668
$ws.Range("H50").Value = 17381
$ws.Range("J50").Value = 24988.334
$ws.Range("L50").Value = 24988.334
$ws.Range("N50").Value = -26238.334
$ws.Range("H59").Value = 58553.777
$ws.Range("J59").Value = 58553.777
$ws.Range("L59").Value = 58553.777
$ws.Range("N59").Value = -60843.777
$ws.Range("H62").Value = 9114.471
$ws.Range("I62").Value = 3300.4546
$ws.Range("J62").Value = 19773.5
$ws.Range("K62").Value = 3300.4546
$ws.Range("L62").Value = 19773.5
$ws.Range("M62").Value = -2676.4546
$ws.Range("N62").Value = -21021.5
$ws.Range("H65").Value = 9114.471
$ws.Range("I65").Value = 3300.4546
$ws.Range("J65").Value = 19773.5
$ws.Range("K65").Value = 16502.273
$ws.Range("L65").Value = 98867.5
$ws.Range("M65").Value = -13382.273
$ws.Range("N65").Value = -105107.5
$ws.Range("H107").Value = 771.619
$ws.Range("I107").Value = 245.25
$ws.Range("K107").Value = 245.25
$ws.Range("M107").Value = 1674.75
$ws.Range("H134").Value = 9587.5
$ws.Range("I134").Value = 9018.096
$ws.Range("K134").Value = 27054.288
$ws.Range("M134").Value = -24519.288

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3636.4
$ws.Range("I75").Value = 3720.5
$ws.Range("K75").Value = 11161.5
$ws.Range("M75").Value = -10163.5
$ws.Range("H78").Value = 3636.4
$ws.Range("I78").Value = 3720.5
$ws.Range("K78").Value = 33484.5
$ws.Range("M78").Value = -28492.5
$ws.Range("H131").Value = 3718.4138
$ws.Range("I131").Value = 2138.75
$ws.Range("J131").Value = 4833.4707
$ws.Range("K131").Value = 6416.25
$ws.Range("L131").Value = 14500.4121
$ws.Range("M131").Value = -1376.25
$ws.Range("N131").Value = -24580.4121

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7850.0557
$ws.Range("I70").Value = 7992.769
$ws.Range("J70").Value = 7479
$ws.Range("K70").Value = 7992.769
$ws.Range("L70").Value = 7479
$ws.Range("M70").Value = -7722.769
$ws.Range("N70").Value = -8019
$ws.Range("H73").Value = 7850.0557
$ws.Range("I73").Value = 7992.769
$ws.Range("J73").Value = 7479
$ws.Range("K73").Value = 7992.769
$ws.Range("L73").Value = 7479
$ws.Range("M73").Value = -7056.769
$ws.Range("N73").Value = -9351
$ws.Range("H97").Value = 1231.8182
$ws.Range("I97").Value = 488.7143
$ws.Range("K97").Value = 488.7143
$ws.Range("M97").Value = 7.28570000000002
$ws.Range("H102").Value = 6984.231
$ws.Range("I102").Value = 11074.5
$ws.Range("J102").Value = 5166.3335
$ws.Range("K102").Value = 11074.5
$ws.Range("L102").Value = 5166.3335
$ws.Range("M102").Value = -9452.5
$ws.Range("N102").Value = -8410.333500000001
$ws.Range("H113").Value = 2836.3333
$ws.Range("J113").Value = 2500
$ws.Range("L113").Value = 2500
$ws.Range("N113").Value = -6840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1131.2
$ws.Range("I93").Value = 920.6923
$ws.Range("K93").Value = 920.6923
$ws.Range("M93").Value = 327.3077
$ws.Range("H100").Value = 1784.8
$ws.Range("I100").Value = 1621.2858
$ws.Range("K100").Value = 1621.2858
$ws.Range("M100").Value = -1080.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7331.909
$ws.Range("I62").Value = 6707.857
$ws.Range("K62").Value = 6707.857
$ws.Range("M62").Value = -6083.857
$ws.Range("H65").Value = 7331.909
$ws.Range("I65").Value = 6707.857
$ws.Range("K65").Value = 33539.285
$ws.Range("M65").Value = -30419.285
$ws.Range("H96").Value = 4226.077
$ws.Range("I96").Value = 3573.1428
$ws.Range("K96").Value = 3573.1428
$ws.Range("M96").Value = -2200.1428
$ws.Range("H107").Value = 927.13336
$ws.Range("I107").Value = 961.1111
$ws.Range("J107").Value = 621.3333
$ws.Range("K107").Value = 2883.3333
$ws.Range("L107").Value = 1863.9999
$ws.Range("M107").Value = -963.3332999999998
$ws.Range("N107").Value = -5703.9999
